# Applies weighting updates to the ranking methodology workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Applicant selectivity (rows 10-12)
$ws.Range("D10").Value = 0.05
$ws.Range("D11").Value = 0.01
$ws.Range("D12").Value = 0.01

# Care quality (rows 13-19)
$ws.Range("D13").Value = 0.05
$ws.Range("D14").Value = 0.01
$ws.Range("D15").Value = 0.01
$ws.Range("D16").Value = 0.01
$ws.Range("D17").Value = 0.01
$ws.Range("D18").Value = 0.01
$ws.Range("D19").Value = 0.01

# Research quality (row 22)
$ws.Range("D22").Value = 0.1

# Update the active cell selection to D22
$ws.Range("D22").Select()
